# Fix bug in dataframe_to_AGS4(): add LLPL worksheet as a regression test case.
$wb = $excel.ActiveWorkbook

# --- Add the new "LLPL" worksheet after the last existing sheet (LOCA) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "LLPL"

# --- Row 1: headers ---
$ws.Range("A1").Value = "HEADING"
$ws.Range("B1").Value = "LOCA_ID"
$ws.Range("C1").Value = "SAMP_TOP"
$ws.Range("D1").Value = "SAMP_REF"
$ws.Range("E1").Value = "SAMP_TYPE"
$ws.Range("F1").Value = "SAMP_ID"
$ws.Range("G1").Value = "SPEC_REF"
$ws.Range("H1").Value = "SPEC_DPTH"
$ws.Range("I1").Value = "LLPL_LL"
$ws.Range("J1").Value = "LLPL_PL"
$ws.Range("K1").Value = "LLPL_PI"
$ws.Range("L1").Value = "LLPL_425"

# --- Row 2: units ---
$ws.Range("A2").Value = "UNIT"
$ws.Range("C2").Value = "m"
$ws.Range("H2").Value = "m"
$ws.Range("I2").Value = "%"
$ws.Range("J2").Value = "%"
$ws.Range("L2").Value = "%"

# --- Row 3: data types ---
$ws.Range("A3").Value = "TYPE"
$ws.Range("B3").Value = "ID"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
$ws.Range("E3").Value = "PA"
$ws.Range("F3").Value = "ID"
$ws.Range("G3").Value = "X"
$ws.Range("H3").Value = "2DP"
$ws.Range("I3").Value = "2SF"
$ws.Range("J3").Value = "XN"
$ws.Range("K3").Value = "2SF"

# --- Row 4 & 5: data ---
$ws.Range("A4").Value = "DATA"
$ws.Range("B4").Value = "327-16A"
$ws.Range("C4").Value = 15.01
$ws.Range("C4").NumberFormat = "0.00"
$ws.Range("E4").Value = "U"
$ws.Range("H4").Value = 15.019
$ws.Range("L4").Value = 15.1234

$ws.Range("A5").Value = "DATA"
$ws.Range("B5").Value = "327-16A"
$ws.Range("C5").Value = 15.14
$ws.Range("C5").NumberFormat = "0.00"
$ws.Range("E5").Value = "U"
$ws.Range("H5").Value = 15.1432
$ws.Range("L5").Value = 15

# --- Cells whose text looks numeric ("15", "1", "2", "45", "25", "20", "40")
# need to be entered as genuine text (not coerced to numbers). Build them
# via a TEXT() formula in a scratch cell, then paste-special as values so
# the resulting cell keeps the default (unformatted) style.
function Set-TextValue($cell, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '=TEXT("' + $text + '","@")'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-TextValue $ws.Range("D4") "15"
Set-TextValue $ws.Range("D5") "15"
Set-TextValue $ws.Range("G4") "1"
Set-TextValue $ws.Range("G5") "2"
Set-TextValue $ws.Range("I4") "45"
Set-TextValue $ws.Range("I5") "40"
Set-TextValue $ws.Range("J4") "25"
Set-TextValue $ws.Range("J5") "20"
Set-TextValue $ws.Range("K4") "20"
Set-TextValue $ws.Range("K5") "20"

# --- Make LLPL the active sheet/tab (activeTab moves from LOCA to LLPL) ---
$ws.Activate()
